# Fill in the Beta pre/post sd, hosp scale and death scale columns
# (B:E) for every state row on Sheet1 that was still missing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @(0.81, 0.009, 0.9, 1.1)          # AL
    3  = @(0.81, 0.16, 1.2, 1.3)           # CO
    4  = @(0.69, 0.07, 1.01, 1)            # FL
    5  = @(0.67, 0.008, 1.2, 1.2)          # GA
    6  = @(0.75, 0.25, 0.1, 1.2)           # ID
    7  = @(0.77, 0.04, 0.8, 1.1)           # KS
    8  = @(0.75, 0.2, 1.5, 1.1)            # KY
    9  = @(1, 0.27, 0.9, 1.6)              # MA
    10 = @(0.66, 0.18, 1.2, 1.1)           # MD
    11 = @(0.7, 0.000001, 0.7, 1)          # ME
    12 = @(0.75, 0.15, 1.2, 1.4)           # MN
    13 = @(0.81, 0.01, 1.15, 1.15)         # MS
    14 = @(0.75, 0.001, 0.7, 0.8)          # MT
    16 = @(0.8, 0.001, 0.6, 1)             # TN
    17 = @(0.5, 0.2, 0.9, 1.1)             # NH
    18 = @(0.66, 0.28, 1, 1)               # NM
    19 = @(0.9, 0.23, 1.2, 1.2)            # OH
    20 = @(0.8, 0.23, 1.3, 1.7)            # OK
    21 = @(0.6, 0.2, 1.4, 1)               # OR
    22 = @(0.54, 0.288, 0.9, 1.4)          # RI
    23 = @(0.57, 0.0001, 1, 1)             # SC
    24 = @(0.86, 0.03, 1.3, 1.3)           # WI
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}

# Restore the active selection used by the author when the workbook was saved.
$ws.Range("J11").Select()
